# HP: ExcelToStringArray extended with write methods
#
# - Rename Sheet1 -> Zipfiles, Sheet2 -> Results
# - Make Results the active sheet/tab
# - Build a "Results" header block on the Results sheet (title + column
#   headers, with a thin bottom border under the header row)
# - Update the remembered selections on both sheets

$wb = $excel.ActiveWorkbook

$wsZip = $wb.Worksheets.Item("Sheet1")
$wsRes = $wb.Worksheets.Item("Sheet2")

# --- Populate the Results sheet -------------------------------------------
# Order matters: it controls the order new entries land in the shared
# string table (Points, Feedback, then Results).
$wsRes.Range("A10").Value = "Order"
$wsRes.Range("B10").Value = "Surname"
$wsRes.Range("C10").Value = "First Name"
$wsRes.Range("D10").Value = "Nr"
$wsRes.Range("F10").Value = "Points"
$wsRes.Range("G10").Value = "Feedback"
$wsRes.Range("B2").Value = "Results"

# Thin bottom border under the whole header row (also touches the blank
# E10 cell so it gets pulled into the used range with the same style).
$wsRes.Range("A10:G10").Borders.Item(9).LineStyle = 1

# --- Rename sheets ----------------------------------------------------------
$wsZip.Name = "Zipfiles"
$wsRes.Name = "Results"

# --- Selections --------------------------------------------------------------
$wsZip.Range("A10:F10").Select()
$wsRes.Range("B3").Select()

# --- Active sheet/tab --------------------------------------------------------
$wsRes.Activate()
